$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4001243.8
$ws.Range("J17").Value = 4001243.8
$ws.Range("L17").Value = 12003731.4
$ws.Range("N17").Value = -12004067.4
$ws.Range("H38").Value = 3127.5
$ws.Range("I38").Value = 182.14285
$ws.Range("K38").Value = 546.4285500000001
$ws.Range("M38").Value = -174.4285500000001
$ws.Range("H39").Value = 309.42856
$ws.Range("I39").Value = 309.42856
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 928.28568
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -632.28568
$ws.Range("N39").ClearContents()
$ws.Range("H42").Value = 226.6
$ws.Range("J42").Value = 45
$ws.Range("L42").Value = 135
$ws.Range("N42").Value = -595
$ws.Range("H97").Value = 1821.5
$ws.Range("J97").Value = 1821.5
$ws.Range("L97").Value = 5464.5
$ws.Range("N97").Value = -6456.5
$ws.Range("H98").Value = 1108.0358
$ws.Range("I98").Value = 975.1667
$ws.Range("K98").Value = 975.1667
$ws.Range("M98").Value = 522.8333
$ws.Range("H99").Value = 207
$ws.Range("I99").Value = 164
$ws.Range("K99").Value = 492
$ws.Range("M99").Value = 1006
$ws.Range("H101").Value = 5153.8
$ws.Range("I101").Value = 2559.3333
$ws.Range("K101").Value = 7677.999899999999
$ws.Range("M101").Value = -6055.999899999999
$ws.Range("H107").Value = 1010.25
$ws.Range("I107").Value = 1047.9474
$ws.Range("K107").Value = 1047.9474
$ws.Range("M107").Value = 872.0526
$ws.Range("H115").Value = 1570.75
$ws.Range("J115").Value = 5000
$ws.Range("L115").Value = 15000
$ws.Range("N115").Value = -18134
$ws.Range("H118").Value = 800
$ws.Range("I118").Value = 200
$ws.Range("K118").Value = 600
$ws.Range("M118").Value = 1057
$ws.Range("H122").Value = 1108.0358
$ws.Range("I122").Value = 975.1667
$ws.Range("K122").Value = 2925.5001
$ws.Range("M122").Value = -475.5001000000002
$ws.Range("H127").Value = 5948.5
$ws.Range("I127").Value = 10000
$ws.Range("K127").Value = 30000
$ws.Range("M127").Value = -25040
$ws.Range("H129").Value = 1311.5
$ws.Range("I129").Value = 1311.5
$ws.Range("K129").Value = 3934.5
$ws.Range("M129").Value = 1065.5
$ws.Range("H132").Value = 1962.5927
$ws.Range("I132").Value = 1730
$ws.Range("J132").Value = 3300
$ws.Range("K132").Value = 5190
$ws.Range("L132").Value = 9900
$ws.Range("M132").Value = -2660
$ws.Range("N132").Value = -14960
$ws.Range("H135").Value = 1145.875
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2120.7368
$ws.Range("I137").Value = 1912.9333
$ws.Range("J137").Value = 2900
$ws.Range("K137").Value = 5738.7999
$ws.Range("L137").Value = 8700
$ws.Range("M137").Value = -3188.7999
$ws.Range("N137").Value = -13800
$ws.Range("H138").Value = 2594.3257
$ws.Range("I138").Value = 2227.923
$ws.Range("J138").Value = 2753.1
$ws.Range("K138").Value = 6683.768999999999
$ws.Range("L138").Value = 8259.299999999999
$ws.Range("M138").Value = -1543.768999999999
$ws.Range("N138").Value = -18539.3
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 1871.4445
$ws.Range("I141").Value = 1871.4445
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5614.333500000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -434.3335000000006
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12185.096
$ws.Range("I61").Value = 11581.588
$ws.Range("K61").Value = 11581.588
$ws.Range("M61").Value = -11369.588
$ws.Range("H74").Value = 6547.654
$ws.Range("I74").Value = 5419.7896
$ws.Range("K74").Value = 5419.7896
$ws.Range("M74").Value = -4545.7896
$ws.Range("H77").Value = 6547.654
$ws.Range("I77").Value = 5419.7896
$ws.Range("K77").Value = 27098.948
$ws.Range("M77").Value = -22730.948
$ws.Range("H96").Value = 93500
$ws.Range("J96").Value = 93500
$ws.Range("L96").Value = 93500
$ws.Range("N96").Value = -98992
$ws.Range("H132").Value = 3134.7368
$ws.Range("I132").Value = 2753.5293
$ws.Range("K132").Value = 8260.5879
$ws.Range("M132").Value = -5730.5879
$ws.Range("H135").Value = 166499.75
$ws.Range("J135").Value = 166499.75
$ws.Range("L135").Value = 166499.75
$ws.Range("N135").Value = -176639.75
$ws.Range("H136").Value = 12185.096
$ws.Range("I136").Value = 11581.588
$ws.Range("K136").Value = 34744.764
$ws.Range("M136").Value = -32194.764

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69159.5
$ws.Range("J2").Value = 69159.5
$ws.Range("L2").Value = 69159.5
$ws.Range("N2").Value = -69385.5
$ws.Range("H44").Value = 10483.333
$ws.Range("J44").Value = 12725
$ws.Range("L44").Value = 12725
$ws.Range("N44").Value = -13719
$ws.Range("H94").Value = 1408.4828
$ws.Range("J94").Value = 3617.25
$ws.Range("L94").Value = 3617.25
$ws.Range("N94").Value = -4519.25
$ws.Range("H99").Value = 4374.409
$ws.Range("J99").Value = 5718.7
$ws.Range("L99").Value = 5718.7
$ws.Range("N99").Value = -8714.700000000001
$ws.Range("H100").Value = 12446.333
$ws.Range("J100").Value = 12446.333
$ws.Range("L100").Value = 12446.333
$ws.Range("N100").Value = -14610.333
$ws.Range("H105").Value = 1370.3077
$ws.Range("I105").Value = 1419.75
$ws.Range("J105").Value = 777
$ws.Range("K105").Value = 1419.75
$ws.Range("L105").Value = 777
$ws.Range("M105").Value = 327.25
$ws.Range("N105").Value = -4271
$ws.Range("H134").Value = 5681.5
$ws.Range("I134").Value = 5113.4165
$ws.Range("K134").Value = 15340.2495
$ws.Range("M134").Value = -12805.2495
$ws.Range("H139").Value = 9999
$ws.Range("J139").Value = 9999
$ws.Range("L139").Value = 9999
$ws.Range("N139").Value = -20279

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4143.2666
$ws.Range("J16").Value = 7797.8
$ws.Range("L16").Value = 7797.8
$ws.Range("N16").Value = -8371.799999999999
$ws.Range("H22").Value = 510.9091
$ws.Range("I22").Value = 273
$ws.Range("J22").Value = 796.4
$ws.Range("K22").Value = 273
$ws.Range("L22").Value = 796.4
$ws.Range("M22").Value = 77
$ws.Range("N22").Value = -1496.4
$ws.Range("H31").Value = 4600.1904
$ws.Range("I31").Value = 4586.625
$ws.Range("J31").Value = 4608.5386
$ws.Range("K31").Value = 4586.625
$ws.Range("L31").Value = 4608.5386
$ws.Range("M31").Value = -4291.625
$ws.Range("N31").Value = -5198.5386
$ws.Range("H34").Value = 4600.1904
$ws.Range("I34").Value = 4586.625
$ws.Range("J34").Value = 4608.5386
$ws.Range("K34").Value = 4586.625
$ws.Range("L34").Value = 4608.5386
$ws.Range("M34").Value = -4384.625
$ws.Range("N34").Value = -5012.5386
$ws.Range("H35").Value = 4261.2
$ws.Range("I35").Value = 3826.5
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 3826.5
$ws.Range("L35").Value = 6000
$ws.Range("M35").Value = -3532.5
$ws.Range("N35").Value = -6588
$ws.Range("H58").Value = 6290.2
$ws.Range("I58").Value = 5345.143
$ws.Range("K58").Value = 5345.143
$ws.Range("M58").Value = -5142.143
$ws.Range("H74").Value = 41767.5
$ws.Range("J74").Value = 41767.5
$ws.Range("L74").Value = 41767.5
$ws.Range("N74").Value = -43515.5
$ws.Range("H77").Value = 41767.5
$ws.Range("J77").Value = 41767.5
$ws.Range("L77").Value = 125302.5
$ws.Range("N77").Value = -134038.5
$ws.Range("H80").Value = 39865.152
$ws.Range("J80").Value = 39865.152
$ws.Range("L80").Value = 39865.152
$ws.Range("N80").Value = -42111.152
$ws.Range("H83").Value = 39865.152
$ws.Range("J83").Value = 39865.152
$ws.Range("L83").Value = 119595.456
$ws.Range("N83").Value = -130827.456
$ws.Range("H99").Value = 4326.5
$ws.Range("J99").Value = 3977.3333
$ws.Range("L99").Value = 3977.3333
$ws.Range("N99").Value = -6973.3333
$ws.Range("H105").Value = 1794.7059
$ws.Range("I105").Value = 1373.8462
$ws.Range("K105").Value = 1373.8462
$ws.Range("M105").Value = 373.1538
$ws.Range("H107").Value = 834.1667
$ws.Range("I107").Value = 801
$ws.Range("K107").Value = 801
$ws.Range("M107").Value = 1119
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 4143.2666
$ws.Range("J113").Value = 7797.8
$ws.Range("L113").Value = 7797.8
$ws.Range("N113").Value = -12137.8
$ws.Range("H122").Value = 4283.9287
$ws.Range("I122").Value = 4921.143
$ws.Range("J122").Value = 3646.7144
$ws.Range("K122").Value = 14763.429
$ws.Range("L122").Value = 10940.1432
$ws.Range("M122").Value = -12313.429
$ws.Range("N122").Value = -15840.1432
$ws.Range("H126").Value = 4326.5
$ws.Range("J126").Value = 3977.3333
$ws.Range("L126").Value = 11931.9999
$ws.Range("N126").Value = -16871.9999
$ws.Range("H132").Value = 5616.1333
$ws.Range("I132").Value = 5662.357
$ws.Range("K132").Value = 16987.071
$ws.Range("M132").Value = -14457.071
$ws.Range("H134").Value = 5909.9473
$ws.Range("I134").Value = 5545.3076
$ws.Range("K134").Value = 16635.9228
$ws.Range("M134").Value = -14100.9228
$ws.Range("H136").Value = 6290.2
$ws.Range("I136").Value = 5345.143
$ws.Range("K136").Value = 16035.429
$ws.Range("M136").Value = -13485.429
$ws.Range("H138").Value = 141934.6
$ws.Range("I138").Value = 76489
$ws.Range("K138").Value = 76489
$ws.Range("M138").Value = -71349

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3999
$ws.Range("J48").Value = 3999
$ws.Range("L48").Value = 11997
$ws.Range("N48").Value = -12497
$ws.Range("H54").Value = 83866.664
$ws.Range("I54").Value = 800
$ws.Range("K54").Value = 2400
$ws.Range("M54").Value = -1841
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H93").Value = 1761.6666
$ws.Range("I93").Value = 1537.5
$ws.Range("K93").Value = 4612.5
$ws.Range("M93").Value = -2740.5
$ws.Range("H113").Value = 1021.1667
$ws.Range("J113").Value = 1075.8
$ws.Range("L113").Value = 3227.4
$ws.Range("N113").Value = -7567.4
$ws.Range("H126").Value = 14600
$ws.Range("H132").Value = 2012.4
$ws.Range("I132").Value = 1772.5555
$ws.Range("K132").Value = 15952.9995
$ws.Range("M132").Value = -13422.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14781.6
$ws.Range("J70").Value = 16329.3
$ws.Range("L70").Value = 16329.3
$ws.Range("N70").Value = -16869.3
$ws.Range("H73").Value = 14781.6
$ws.Range("J73").Value = 16329.3
$ws.Range("L73").Value = 16329.3
$ws.Range("N73").Value = -18201.3
$ws.Range("H99").Value = 10541
$ws.Range("I99").Value = 9649.200000000001
$ws.Range("K99").Value = 9649.200000000001
$ws.Range("M99").Value = -7403.200000000001
$ws.Range("H113").Value = 503520.75
$ws.Range("J113").Value = 4192.5
$ws.Range("L113").Value = 4192.5
$ws.Range("N113").Value = -8532.5
$ws.Range("H122").Value = 7062.3335
$ws.Range("I122").Value = 7380.4165
$ws.Range("K122").Value = 22141.2495
$ws.Range("M122").Value = -19691.2495
$ws.Range("H132").Value = 1045.0769
$ws.Range("I132").Value = 945.1111
$ws.Range("J132").Value = 1270
$ws.Range("K132").Value = 2835.3333
$ws.Range("L132").Value = 3810
$ws.Range("M132").Value = -305.3332999999998
$ws.Range("N132").Value = -8870
$ws.Range("H137").Value = 61748.75
$ws.Range("I137").Value = 61748.75
$ws.Range("K137").Value = 61748.75
$ws.Range("M137").Value = -56648.75
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2342.7778
$ws.Range("I7").Value = 2342.7778
$ws.Range("K7").Value = 2342.7778
$ws.Range("M7").Value = -2230.7778
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H16").Value = 1808.4706
$ws.Range("I16").Value = 1667.4286
$ws.Range("K16").Value = 1667.4286
$ws.Range("M16").Value = -1497.4286
$ws.Range("H20").Value = 103333
$ws.Range("J20").Value = 103333
$ws.Range("L20").Value = 103333
$ws.Range("N20").Value = -103785
$ws.Range("H40").Value = 2075.2173
$ws.Range("I40").Value = 1982.9048
$ws.Range("J40").Value = 3044.5
$ws.Range("K40").Value = 1982.9048
$ws.Range("L40").Value = 3044.5
$ws.Range("M40").Value = -1846.9048
$ws.Range("N40").Value = -3316.5
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178
$ws.Range("H122").Value = 4295.0625
$ws.Range("I122").Value = 3538
$ws.Range("J122").Value = 4749.3
$ws.Range("K122").Value = 10614
$ws.Range("L122").Value = 14247.9
$ws.Range("M122").Value = -8164
$ws.Range("N122").Value = -19147.9
$ws.Range("H126").Value = 2342.7778
$ws.Range("I126").Value = 2342.7778
$ws.Range("K126").Value = 7028.3334
$ws.Range("M126").Value = -4558.3334
$ws.Range("H132").Value = 9443.729499999999
$ws.Range("I132").Value = 10140.533
$ws.Range("J132").Value = 6457.4287
$ws.Range("K132").Value = 30421.599
$ws.Range("L132").Value = 19372.2861
$ws.Range("M132").Value = -27891.599
$ws.Range("N132").Value = -24432.2861
$ws.Range("H136").Value = 3323.152
$ws.Range("I136").Value = 3099.6
$ws.Range("J136").Value = 4813.5
$ws.Range("K136").Value = 9298.799999999999
$ws.Range("L136").Value = 14440.5
$ws.Range("M136").Value = -6748.799999999999
$ws.Range("N136").Value = -19540.5
$ws.Range("H140").Value = 97246.5
$ws.Range("J140").Value = 97246.5
$ws.Range("L140").Value = 97246.5
$ws.Range("N140").Value = -107606.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5546.3335
$ws.Range("I62").Value = 5125.067
$ws.Range("J62").Value = 7652.6665
$ws.Range("K62").Value = 5125.067
$ws.Range("L62").Value = 7652.6665
$ws.Range("M62").Value = -4501.067
$ws.Range("N62").Value = -8900.666499999999
$ws.Range("H65").Value = 5546.3335
$ws.Range("I65").Value = 5125.067
$ws.Range("J65").Value = 7652.6665
$ws.Range("K65").Value = 25625.335
$ws.Range("L65").Value = 38263.3325
$ws.Range("M65").Value = -22505.335
$ws.Range("N65").Value = -44503.3325
$ws.Range("H107").Value = 2181.1667
$ws.Range("I107").Value = 1495.6666
$ws.Range("J107").Value = 2866.6667
$ws.Range("K107").Value = 4486.9998
$ws.Range("L107").Value = 8600.000100000001
$ws.Range("M107").Value = -2566.9998
$ws.Range("N107").Value = -12440.0001
$ws.Range("H122").Value = 4832.1724
$ws.Range("I122").Value = 4494.5
$ws.Range("K122").Value = 13483.5
$ws.Range("M122").Value = -11033.5
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040
$ws.Range("H132").Value = 5484.394
$ws.Range("I132").Value = 6080.75
$ws.Range("J132").Value = 4566.923
$ws.Range("K132").Value = 18242.25
$ws.Range("L132").Value = 13700.769
$ws.Range("M132").Value = -15712.25
$ws.Range("N132").Value = -18760.769
$ws.Range("H136").Value = 4663.8
$ws.Range("I136").Value = 2714.5417
$ws.Range("J136").Value = 12460.833
$ws.Range("K136").Value = 8143.625100000001
$ws.Range("L136").Value = 37382.499
$ws.Range("M136").Value = -5593.625100000001
$ws.Range("N136").Value = -42482.499
